$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the cell content: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Move the active selection to E8, matching the saved view state
$ws.Activate()
$ws.Range("E8").Select()
